$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fix product list ordering ("1-torta,2-queque," -> "2-queque,1-torta,")
$ws.Range("A2").Value = "2-queque,1-torta,"

# Row 8: correct status from "En Proceso" to "Pendiente"
$ws.Range("G8").Value = "Pendiente"

# Row 8: correct "Precio cancelado" from 10 to 0
$ws.Range("I8").Value = 0.0
